# Smart Meters.xlsx update — "smart meters, rhi, daily demand, chp"
#
# Adds two new monthly data rows (Aug-2023 / 1-Aug-2023 and Sep-2023 / 1-Sep-2023)
# to the bottom of the table, corrects the previously-estimated Jun-2023 (row 81)
# figures now that firmer numbers are available, and leaves the sheet scrolled /
# selected the way the author left it after typing the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 81 (1-Jun-2023) — D81 revised from 1,133,143 to 1,333,143; B81 is the
#    shared formula C81+D81 and recalculates automatically.
# ---------------------------------------------------------------------------
$ws.Range("D81").Value = 1333143

# ---------------------------------------------------------------------------
# 2) Append row 83 (1-Aug-2023) and row 84 (1-Sep-2023). Copy row 82's
#    formatting down (keeps the same cell styles: date format / number
#    formats) then fill in the new values and formulas.
# ---------------------------------------------------------------------------
$ws.Range("A82:D82").Copy()
$ws.Range("A83:D84").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A83").Value = 45139
$ws.Range("B83").Formula = "=C83+D83"
$ws.Range("C83").Value = 357329
$ws.Range("D83").Value = 1359112

$ws.Range("A84").Value = 45170
$ws.Range("B84").Formula = "=C84+D84"
$ws.Range("C84").Value = 362961
$ws.Range("D84").Value = 1372266

# ---------------------------------------------------------------------------
# 3) Leave the view scrolled/selected the way it was left after the edit.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D81").Select()
